$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the shared-string table in the order the strings were first typed:
# "Flower" (Flower power-up rows), then "NONE", then "iup".
$ws.Range("H11").Value = "Flower"
$ws.Range("H3").Value = "NONE"
$ws.Range("J5").Value = "iup"

# Row 1
$ws.Range("G1").Value = 4671
$ws.Range("I1").Value = 4672

# Row 3
$ws.Range("G3").Value = 22200

# Row 4
$ws.Range("G4").Value = 22300
$ws.Range("H4").Value = "NONE"

# Row 5
$ws.Range("G5").Value = 22500
$ws.Range("H5").Value = "NONE"
$ws.Range("I5").Value = 22700

# Row 6
$ws.Range("G6").Value = 22900
$ws.Range("H6").Value = "NONE"
$ws.Range("I6").Value = 22900
$ws.Range("J6").Value = "NONE"

# Row 7
$ws.Range("G7").Value = 23100
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = "NONE"
$ws.Range("H7").Style = "Normal"

# Row 11
$ws.Range("G11").Value = 23500

# Row 12
$ws.Range("G12").Value = 23600
$ws.Range("H12").Value = "Flower"

# Row 13
$ws.Range("G13").Value = 23800
$ws.Range("H13").Value = "Flower"

# Row 14
$ws.Range("G14").Value = 24200
$ws.Range("H14").Value = "Flower"

# Row 15
$ws.Range("G15").Value = 24400
$ws.Range("H15").Value = "Flower"

# Update the active selection to match the edited workbook
[void]$ws.Range("J7").Select()
